# Update odds values on Sheet1 to match the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.45
$ws.Range("H2").Value = 2.82
$ws.Range("K2").Value = 1.77
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.63
$ws.Range("T2").Value = 1.33
$ws.Range("U2").Value = 5.4
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 1.1

# Row 3
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 2.85
$ws.Range("K3").Value = 1.87
$ws.Range("T3").Value = 1.41

# Row 4
$ws.Range("K4").Value = 1.8

# Row 5
$ws.Range("K5").Value = 1.8

# Row 6
$ws.Range("J6").Value = 1.92

# Row 7
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 4.33
$ws.Range("J7").Value = 2.37
$ws.Range("K7").Value = 2.3
$ws.Range("L7").Value = 4.5
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 4
$ws.Range("S7").Value = 1.73
$ws.Range("T7").Value = 2.1
$ws.Range("U7").Value = 2.05
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 3
$ws.Range("X7").Value = 1.4
$ws.Range("Y7").Value = 1.33
$ws.Range("Z7").Value = 3.25
$ws.Range("AA7").Value = 1.7
$ws.Range("AB7").Value = 2.05
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 9
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 13
$ws.Range("AJ7").Value = 7.5
$ws.Range("AM7").Value = 13
$ws.Range("AN7").Value = 23
$ws.Range("AQ7").Value = 34
$ws.Range("AR7").Value = 34
$ws.Range("AS7").Value = 201

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 2.87

# Row 9
$ws.Range("G9").Value = 1.52
$ws.Range("H9").Value = 4
$ws.Range("K9").Value = 2.37
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("Y9").Value = 1.33
$ws.Range("Z9").Value = 3.25
$ws.Range("AF9").Value = 12

# Row 13
$ws.Range("I13").Value = 2.05

# Row 14
$ws.Range("I14").Value = 2.01
$ws.Range("S14").Value = 1.77
$ws.Range("T14").Value = 1.97
$ws.Range("W14").Value = 3
$ws.Range("X14").Value = 1.36

# Row 15
$ws.Range("I15").Value = 1.69
$ws.Range("AS15").Value = 1000

# Row 16
$ws.Range("Q16").Value = 1.8
$ws.Range("S16").Value = 2.4
$ws.Range("T16").Value = 1.53
$ws.Range("W16").Value = 4.5
$ws.Range("X16").Value = 1.18

# Row 19
$ws.Range("G19").Value = 2.57
$ws.Range("H19").Value = 2.82
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 3.4
$ws.Range("K19").Value = 1.95
$ws.Range("L19").Value = 3.75
$ws.Range("M19").Value = 1.11
$ws.Range("N19").Value = 6.5
$ws.Range("Q19").Value = 1.82
$ws.Range("R19").Value = 1.92
$ws.Range("AD19").Value = 11
$ws.Range("AQ19").Value = 29
